# The "Authors" column (E) for rows 2-12 was re-generated by the new
# elastic-search-backed author lookup. The new values are identical to the
# previous ones except that every comma-separated author entry now has two
# extra spaces of padding before it (a side effect of the new formatter).
# Reproduce that by taking the current value of each cell and inserting two
# extra spaces after every "," + whitespace run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Range("E$r")
    $orig = $cell.Value2
    $updated = $orig -replace ',(\s+)', ',  $1'
    $cell.Value = $updated
}
